$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(5610, 45947),
    @(5565, 45947.01041666666),
    @(5518, 45947.02083333334),
    @(5490, 45947.03125),
    @(5449, 45947.04166666666),
    @(5455, 45947.05208333334),
    @(5399, 45947.0625),
    @(5404, 45947.07291666666),
    @(5415, 45947.08333333334),
    @(5425, 45947.09375),
    @(5412, 45947.10416666666),
    @(5365, 45947.11458333334),
    @(5410, 45947.125),
    @(5435, 45947.13541666666),
    @(5441, 45947.14583333334),
    @(5500, 45947.15625),
    @(5529, 45947.16666666666),
    @(5641, 45947.17708333334),
    @(5679, 45947.1875),
    @(5790, 45947.19791666666),
    @(5990, 45947.20833333334),
    @(6166, 45947.21875),
    @(6344, 45947.22916666666),
    @(6495, 45947.23958333334),
    @(6808, 45947.25),
    @(6949, 45947.26041666666),
    @(7077, 45947.27083333334),
    @(7105, 45947.28125),
    @(7248, 45947.29166666666),
    @(7300, 45947.30208333334),
    @(7269, 45947.3125),
    @(7098, 45947.32291666666),
    @(7066, 45947.33333333334),
    @(6988, 45947.34375),
    @(6917, 45947.35416666666),
    @(6797, 45947.36458333334),
    @(6672, 45947.375),
    @(6473, 45947.38541666666),
    @(6342, 45947.39583333334),
    @(6207, 45947.40625),
    @(5994, 45947.41666666666)
)

$startRow = 2
$endRow = $startRow + $data.Count - 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Make sure the extended range for column B keeps the same timestamp number format
$fmt = $ws.Range("B2").NumberFormat
$ws.Range("B2:B" + $endRow).NumberFormat = $fmt
